$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.008.13"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.85%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.673.58"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.96%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.532"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.83%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  +2.83%  "
$ws.Range("E9").Value = "  +1.80%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.20"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.08%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0890"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.17%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.910.50"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.05%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.671.29"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.94%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.09"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.05%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "65.79"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.81%  "
$ws.Range("E16").Value = "  +1.85%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.032.14"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.93%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "233.32"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.61%  "
$ws.Range("E19").Value = "  +1.56%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.78"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.37%  "
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.46"
$ws.Range("D22").Style = "Normal"
$ws.Range("E23").Value = "  +0.42%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.73"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.19%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.15"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.03%  "
$ws.Range("B27").Value = "Stellar"
$ws.Range("C27").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.116"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.86"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.47%  "
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("E30").Value = "  +0.91%  "
$ws.Range("E31").Value = "  +1.64%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.32"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.78%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.452.36"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.59%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.16"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.17%  "
$ws.Range("E35").Value = "  +5.42%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.40"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.44%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.898"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.35%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.565"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.22%  "
$ws.Range("E39").Value = "  +1.63%  "
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.30"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.91%  "
$ws.Range("E43").Value = "  +7.63%  "
$ws.Range("E44").Value = "  +5.24%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.816.26"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.88%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.783"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.79%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.61"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.65%  "
$ws.Range("E48").Value = "  +1.21%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.100"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.30%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0507"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.60"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.35%  "
